$d = $word.ActiveDocument

# Walk every section's headers and footers and rename the inline logo
# pictures found there:
#   - the Pearson Edexcel logo (currently "image2.png")  -> "image1.png"
#   - the BTEC logo            (currently "image1.jpg")  -> "image2.jpg"
# Matching is done on AlternativeText (the picture's description), which is
# stable regardless of which physical header/footer part backs a given
# Headers/Footers collection index.

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hf = $sec.Headers.Item($i)
        if ($hf.Exists) {
            for ($j = 1; $j -le $hf.Range.InlineShapes.Count; $j++) {
                $shp = $hf.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
                elseif ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $hf = $sec.Footers.Item($i)
        if ($hf.Exists) {
            for ($j = 1; $j -le $hf.Range.InlineShapes.Count; $j++) {
                $shp = $hf.Range.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
                elseif ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
